$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 (index A8:L8) - new job entry
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "UAT"
$ws.Range("C8").Value = "TRD"
$ws.Range("D8").Value = "UAT"
$ws.Range("E8").Value = "TRD"
$ws.Range("F8").Value = "Job_Details"
$ws.Range("G8").Value = "Partition Copy"
$ws.Range("H8").Value = "28-04-2020 09:36:03"
$ws.Range("I8").Value = "28-04-2020 09:36:03"
$ws.Range("J8").Value = "In Progress"
$ws.Range("K8").Value = "Export in Progress"
$ws.Range("L8").Value = "Import in Progress"

# Row 9 (index A9:L9) - new job entry
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "INT"
$ws.Range("C9").Value = "COVID"
$ws.Range("D9").Value = "UAT"
$ws.Range("E9").Value = "HACK"
$ws.Range("F9").Value = "Transactions"
$ws.Range("G9").Value = "Partition Copy"
$ws.Range("H9").Value = "28-04-2020 10:05:55"
$ws.Range("I9").Value = "28-04-2020 10:05:55"
$ws.Range("J9").Value = "In Progress"
$ws.Range("K9").Value = "Export in Progress"
$ws.Range("L9").Value = "Import in Progress"
